$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.987.78'
$ws.Range('E2').Value = '  -0.52%  '
$ws.Range('D3').Value = '2.607.56'
$ws.Range('E3').Value = '  -1.38%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Formula = '''590.33'
$ws.Range('E5').Value = '  -1.38%  '
$ws.Range('D6').Formula = '''165.11'
$ws.Range('E6').Value = '  -0.34%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  -2.14%  '
$ws.Range('D9').Value = '2.607.84'
$ws.Range('E9').Value = '  -1.35%  '
$ws.Range('E10').Value = '  -5.15%  '
$ws.Range('E11').Value = '  +1.45%  '
$ws.Range('D12').Formula = '''0.361'
$ws.Range('E12').Value = '  -0.61%  '
$ws.Range('D13').Formula = '''5.18'
$ws.Range('E13').Value = '  -0.72%  '
$ws.Range('D14').Formula = '''27.23'
$ws.Range('E14').Value = '  -2.60%  '
$ws.Range('D15').Value = '3.080.57'
$ws.Range('E15').Value = '  -1.37%  '
$ws.Range('E16').Value = '  -2.89%  '
$ws.Range('D17').Value = '67.052.51'
$ws.Range('E17').Value = '  -0.41%  '
$ws.Range('D18').Value = '2.607.58'
$ws.Range('E18').Value = '  -0.87%  '
$ws.Range('E19').Value = '  -1.03%  '
$ws.Range('D20').Formula = '''7.79'
$ws.Range('E20').Value = '  -0.85%  '
$ws.Range('D21').Formula = '''353.42'
$ws.Range('E21').Value = '  -2.80%  '
$ws.Range('D22').Formula = '''4.26'
$ws.Range('E22').Value = '  -3.01%  '
$ws.Range('E23').Value = '  -3.68%  '
$ws.Range('D24').Formula = '''10.52'
$ws.Range('E24').Value = '  -5.16%  '
$ws.Range('E25').Value = '  -0.11%  '
$ws.Range('E26').Value = '  -4.65%  '
$ws.Range('D27').Formula = '''68.89'
$ws.Range('E27').Value = '  -2.76%  '
$ws.Range('D28').Value = '2.748.21'
$ws.Range('E28').Value = '  -1.02%  '
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('D30').Value = '0.0₃0993'
$ws.Range('E30').Value = '  -2.93%  '
$ws.Range('D31').Formula = '''539.06'
$ws.Range('E31').Value = '  -2.67%  '
$ws.Range('D32').Formula = '''7.85'
$ws.Range('E32').Value = '  -2.20%  '
$ws.Range('E33').Value = '  -3.74%  '
$ws.Range('E34').Value = '  -2.71%  '
$ws.Range('D35').Formula = '''0.133'
$ws.Range('E35').Value = '  +0.48%  '
$ws.Range('D36').Formula = '''0.999'
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('E37').Value = '  -3.73%  '
$ws.Range('D38').Formula = '''157.39'
$ws.Range('E38').Value = '  -0.17%  '
$ws.Range('D39').Formula = '''18.87'
$ws.Range('E39').Value = '  -2.86%  '
$ws.Range('E40').Value = '  -2.29%  '
$ws.Range('E41').Value = '  +1.69%  '
$ws.Range('E42').Value = '  -1.24%  '
$ws.Range('D43').Formula = '''5.11'
$ws.Range('E43').Value = '  -2.78%  '
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('E45').Value = '  -4.87%  '
$ws.Range('D46').Value = '0.0₆0298'
$ws.Range('E46').Value = '  -1.41%  '
$ws.Range('D47').Formula = '''150.81'
$ws.Range('E47').Value = '  -2.15%  '
$ws.Range('E48').Value = '  -3.95%  '
$ws.Range('E49').Value = '  -3.12%  '
$ws.Range('E50').Value = '  -2.19%  '
$ws.Range('D51').Formula = '''0.0769'
$ws.Range('E51').Value = '  -1.10%  '
